$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Source"
$ws.Range("B1").Value = "C"
$ws.Range("C1").Value = "FFR"
$ws.Range("D1").Value = "LF"

# Row 2: C Lag
$ws.Range("A2").Value = "C Lag"
$ws.Range("B2").Value = "-0.46***"
$ws.Range("C2").Value = "3.79"
$ws.Range("D2").Value = "-6.09"

# Row 3: FFR Lag
$ws.Range("A3").Value = "FFR Lag"
$ws.Range("B3").Value = "-0.01"
$ws.Range("C3").Value = "1.6***"
$ws.Range("D3").Value = "0.5***"

# Row 4: LF Lag
$ws.Range("A4").Value = "LF Lag"
$ws.Range("B4").Value = "0.04*"
$ws.Range("C4").Value = "3.53*"
$ws.Range("D4").Value = "0.54*"
